$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 295.72
$ws.Range("I53").Value = 250.18182
$ws.Range("J53").Value = 331.5
$ws.Range("K53").Value = 250.18182
$ws.Range("L53").Value = 331.5
$ws.Range("M53").Value = 386.81818
$ws.Range("N53").Value = -1605.5
$ws.Range("H132").Value = 4467783
$ws.Range("I132").Value = 5105307.5
$ws.Range("J132").Value = 5113
$ws.Range("K132").Value = 15315922.5
$ws.Range("L132").Value = 15339
$ws.Range("M132").Value = -15313392.5
$ws.Range("N132").Value = -20399
$ws.Range("H137").Value = 1577
$ws.Range("I137").Value = 1573.8125
$ws.Range("K137").Value = 4721.4375
$ws.Range("M137").Value = -2171.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25416.195
$ws.Range("I32").Value = 4237.9355
$ws.Range("K32").Value = 4237.9355
$ws.Range("M32").Value = -3950.9355
$ws.Range("H61").Value = 1467.3055
$ws.Range("I61").Value = 1388.1724
$ws.Range("J61").Value = 1795.1428
$ws.Range("K61").Value = 1388.1724
$ws.Range("L61").Value = 1795.1428
$ws.Range("M61").Value = -1176.1724
$ws.Range("N61").Value = -2219.1428
$ws.Range("H74").Value = 3412.647
$ws.Range("I74").Value = 1774.5454
$ws.Range("J74").Value = 6415.8335
$ws.Range("K74").Value = 1774.5454
$ws.Range("L74").Value = 6415.8335
$ws.Range("M74").Value = -900.5454
$ws.Range("N74").Value = -8163.8335
$ws.Range("H77").Value = 3412.647
$ws.Range("I77").Value = 1774.5454
$ws.Range("J77").Value = 6415.8335
$ws.Range("K77").Value = 8872.726999999999
$ws.Range("L77").Value = 32079.1675
$ws.Range("M77").Value = -4504.726999999999
$ws.Range("N77").Value = -40815.1675
$ws.Range("H132").Value = 2250.2114
$ws.Range("I132").Value = 2090.7646
$ws.Range("J132").Value = 2551.389
$ws.Range("K132").Value = 6272.293799999999
$ws.Range("L132").Value = 7654.167
$ws.Range("M132").Value = -3742.293799999999
$ws.Range("N132").Value = -12714.167
$ws.Range("H136").Value = 1467.3055
$ws.Range("I136").Value = 1388.1724
$ws.Range("J136").Value = 1795.1428
$ws.Range("K136").Value = 4164.5172
$ws.Range("L136").Value = 5385.428400000001
$ws.Range("M136").Value = -1614.5172
$ws.Range("N136").Value = -10485.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1626.75
$ws.Range("I99").Value = 1235.3334
$ws.Range("K99").Value = 1235.3334
$ws.Range("M99").Value = 262.6666
$ws.Range("H107").Value = 125059016
$ws.Range("I107").Value = 200091820
$ws.Range("K107").Value = 200091820
$ws.Range("M107").Value = -200089900
$ws.Range("H134").Value = 2958.0386
$ws.Range("I134").Value = 3057.7827
$ws.Range("J134").Value = 2193.3333
$ws.Range("K134").Value = 9173.348100000001
$ws.Range("L134").Value = 6579.999899999999
$ws.Range("M134").Value = -6638.348100000001
$ws.Range("N134").Value = -11649.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 55706.223
$ws.Range("I31").Value = 43023.707
$ws.Range("J31").Value = 70200.52
$ws.Range("K31").Value = 43023.707
$ws.Range("L31").Value = 70200.52
$ws.Range("M31").Value = -42728.707
$ws.Range("N31").Value = -70790.52
$ws.Range("H34").Value = 55706.223
$ws.Range("I34").Value = 43023.707
$ws.Range("J34").Value = 70200.52
$ws.Range("K34").Value = 43023.707
$ws.Range("L34").Value = 70200.52
$ws.Range("M34").Value = -42821.707
$ws.Range("N34").Value = -70604.52
$ws.Range("H58").Value = 8132.205
$ws.Range("I58").Value = 1373.4546
$ws.Range("K58").Value = 1373.4546
$ws.Range("M58").Value = -1170.4546
$ws.Range("H118").Value = 44964.5
$ws.Range("J118").Value = 44964.5
$ws.Range("L118").Value = 44964.5
$ws.Range("N118").Value = -48278.5
$ws.Range("H132").Value = 4393.5293
$ws.Range("I132").Value = 4665.091
$ws.Range("J132").Value = 3895.6667
$ws.Range("K132").Value = 13995.273
$ws.Range("L132").Value = 11687.0001
$ws.Range("M132").Value = -11465.273
$ws.Range("N132").Value = -16747.0001
$ws.Range("H134").Value = 1101.65
$ws.Range("I134").Value = 1002.6061
$ws.Range("J134").Value = 1568.5714
$ws.Range("K134").Value = 3007.8183
$ws.Range("L134").Value = 4705.7142
$ws.Range("M134").Value = -472.8182999999999
$ws.Range("N134").Value = -9775.7142
$ws.Range("H136").Value = 8132.205
$ws.Range("I136").Value = 1373.4546
$ws.Range("K136").Value = 4120.3638
$ws.Range("M136").Value = -1570.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1400.2858
$ws.Range("I17").Value = 1133.3334
$ws.Range("K17").Value = 3400.0002
$ws.Range("M17").Value = -3231.0002
$ws.Range("H113").Value = 549.13336
$ws.Range("I113").Value = 494.4375
$ws.Range("J113").Value = 611.6429000000001
$ws.Range("K113").Value = 1483.3125
$ws.Range("L113").Value = 1834.9287
$ws.Range("M113").Value = 686.6875
$ws.Range("N113").Value = -6174.9287
$ws.Range("H129").Value = 234574.1
$ws.Range("I129").Value = 6685.1113
$ws.Range("J129").Value = 392343.38
$ws.Range("K129").Value = 20055.3339
$ws.Range("L129").Value = 1177030.14
$ws.Range("M129").Value = -15055.3339
$ws.Range("N129").Value = -1187030.14
$ws.Range("H137").Value = 47884.12
$ws.Range("I137").Value = 116686.664
$ws.Range("J137").Value = 9182.6875
$ws.Range("K137").Value = 350059.992
$ws.Range("L137").Value = 27548.0625
$ws.Range("M137").Value = -344959.992
$ws.Range("N137").Value = -37748.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 999
$ws.Range("J122").Value = 999
$ws.Range("L122").Value = 2997
$ws.Range("N122").Value = -7897
$ws.Range("H126").Value = 2857.6875
$ws.Range("I126").Value = 2870.2856
$ws.Range("J126").Value = 2847.889
$ws.Range("K126").Value = 8610.856800000001
$ws.Range("L126").Value = 8543.667000000001
$ws.Range("M126").Value = -6140.856800000001
$ws.Range("N126").Value = -13483.667
$ws.Range("H132").Value = 3569.45
$ws.Range("I132").Value = 3116
$ws.Range("K132").Value = 9348
$ws.Range("M132").Value = -6818

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3100
$ws.Range("I122").Value = 3780
$ws.Range("K122").Value = 11340
$ws.Range("M122").Value = -8890
$ws.Range("H132").Value = 4017.35
$ws.Range("I132").Value = 4709.857
$ws.Range("J132").Value = 2401.5
$ws.Range("K132").Value = 14129.571
$ws.Range("L132").Value = 7204.5
$ws.Range("M132").Value = -11599.571
$ws.Range("N132").Value = -12264.5
$ws.Range("H136").Value = 1242.6316
$ws.Range("I136").Value = 1042.1666
$ws.Range("J136").Value = 1994.375
$ws.Range("K136").Value = 3126.4998
$ws.Range("L136").Value = 5983.125
$ws.Range("M136").Value = -576.4998000000001
$ws.Range("N136").Value = -11083.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 661.6667
$ws.Range("I113").Value = 448.8889
$ws.Range("J113").Value = 874.44446
$ws.Range("K113").Value = 1346.6667
$ws.Range("L113").Value = 2623.33338
$ws.Range("M113").Value = 823.3333
$ws.Range("N113").Value = -6963.33338
$ws.Range("H132").Value = 1184.8182
$ws.Range("I132").Value = 788
$ws.Range("J132").Value = 2243
$ws.Range("K132").Value = 2364
$ws.Range("L132").Value = 6729
$ws.Range("M132").Value = 166
$ws.Range("N132").Value = -11789
$ws.Range("H136").Value = 1147
$ws.Range("I136").Value = 993.2
$ws.Range("J136").Value = 1403.3334
$ws.Range("K136").Value = 2979.6
$ws.Range("L136").Value = 4210.0002
$ws.Range("M136").Value = -429.6000000000004
$ws.Range("N136").Value = -9310.0002
